$d = $word.ActiveDocument

# Common run-properties fragment (Times New Roman everywhere in this doc)
$rPrXml = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>'
$wDocOpen = "<?xml version='1.0' encoding='UTF-8' standalone='yes'?><pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>"
$wDocClose = "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

# -----------------------------------------------------------------
# 0) Remove the _GoBack bookmark from its old location (end of the
#    last paragraph) first, while its name is still unambiguous -
#    a second "_GoBack" bookmark is (re)created later in step 2.
# -----------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# -----------------------------------------------------------------
# 1) "Steven Qin: ..." paragraph - keep only the "Steven Qin: " lead-in,
#    drop the two trailing runs of commentary text.
# -----------------------------------------------------------------
$stevenPara = $d.Paragraphs.Item(7)
$stevenFull = $stevenPara.Range
$stevenTrim = $d.Range($stevenFull.Start, $stevenFull.End - 1)
$stevenTrim.Text = "Steven Qin: "

# -----------------------------------------------------------------
# 2) Repurpose the (now) "Shawn Yap: " paragraph into the new
#    lesson-learned paragraph: its single run becomes a <w:tab/>,
#    followed by five new text runs and the _GoBack bookmark.
# -----------------------------------------------------------------
$shawnLeadPara = $d.Paragraphs.Item(8)
$shawnLeadFull = $shawnLeadPara.Range
$shawnLeadRange = $d.Range($shawnLeadFull.Start, $shawnLeadFull.End - 1)

$runsXml = ""
$runsXml += "<w:r>$rPrXml<w:tab/></w:r>"
$runsXml += "<w:r>$rPrXml<w:t xml:space='preserve'>Though this project, I was able to get a much deeper understanding on using java, and communicate from one language to another in one program. </w:t></w:r>"
$runsXml += "<w:r>$rPrXml<w:t xml:space='preserve'>Based on the structure of this project, the design was changed several time throughout the process. A good qualified original design can really make the modifying process easier.  </w:t></w:r>"
$runsXml += "<w:r>$rPrXml<w:t xml:space='preserve'>We did a much better job on balancing jobs among all three team members. </w:t></w:r>"
$runsXml += "<w:r>$rPrXml<w:t xml:space='preserve'>Although this is a more difficult task, it was not as stressful as the first one. We learned to do better design and time arrangement. </w:t></w:r>"
$runsXml += "<w:r>$rPrXml<w:t>We are also more familiar with the software we used throughout the second project.</w:t></w:r>"
$runsXml += "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>"

$xmlLesson = "$wDocOpen<w:p>$runsXml</w:p>$wDocClose"
[void]$shawnLeadRange.InsertXML($xmlLesson)

# -----------------------------------------------------------------
# 3) Insert a brand-new "Shawn Yap: " paragraph right after the one
#    that was just repurposed.
# -----------------------------------------------------------------
$lessonPara = $d.Paragraphs.Item(8)
[void]$lessonPara.Range.InsertParagraphAfter()
$newBlankPara = $d.Paragraphs.Item(9)
$newBlankFull = $d.Range($newBlankPara.Range.Start, $newBlankPara.Range.End)

$shawnParaXml = "<w:p><w:pPr>$rPrXml</w:pPr><w:r>$rPrXml<w:t xml:space='preserve'>Shawn Yap: </w:t></w:r></w:p>"
$xmlShawn = "$wDocOpen$shawnParaXml$wDocClose"
[void]$newBlankFull.InsertXML($xmlShawn)
